# Mexico Liga MX workbook update
# - Swaps the full data (all columns except the sequential id in column A) between
#   nine pairs of rows (these rows had their Home/Away data entered in the wrong
#   physical row).
# - Removes two stale placeholder rows (originally rows 288 and 289) that had no
#   result yet.
# - Re-sequences the id column (A) so it stays a contiguous 0..N sequence after
#   the two rows are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap entire row content (columns B:AC) between these row pairs.
$pairs = @(
    @(2,3),
    @(34,35),
    @(72,73),
    @(94,95),
    @(175,176),
    @(193,194),
    @(251,252),
    @(264,265),
    @(272,273)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range("B$r1`:AC$r1")
    $rng2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value = $v2
    $rng2.Value = $v1
}

# 2) Delete the two stale rows (288 and 289); everything below shifts up.
$ws.Rows.Item(289).Delete()
$ws.Rows.Item(288).Delete()

# 3) Re-sequence column A (id) to stay 0..N across the now-shorter sheet.
$lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
